# Update the explanation text about where the MongoDB dump files come from
# on the "Restauração BD" worksheet (cell B4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Restauração BD")

$ws.Range("B4").Value = 'Copiar para a pasta dump do MongoDB o conteúdo da pasta "dump\SWChallenge" do projeto (todos os arquivos .bson e .json)'
